# metadata and cellranger scripts added
# Rework the Visium brain metadata sheet:
#   - replace "condition" (old col C) / "mice_id/bio_origin" (old col D)
#     layout with a new sample_id / section_id / condition layout
#   - drop the now-unused CF2_B / CF7_B / CF1_B / CG9_B / CG8_B strings
#   - clear the Times New Roman override that used to sit on column D
#   - resize columns C and D to fit the new (wider) values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row -----------------------------------------------------
$ws.Cells.Item(1, 1).Value = "sample_no"
$ws.Cells.Item(1, 2).Value = "sample_name"
$ws.Cells.Item(1, 3).Value = "sample_id"
$ws.Cells.Item(1, 4).Value = "section_id"
$ws.Cells.Item(1, 5).Value = "condition"

# ---- data rows --------------------------------------------------------
# columns: A=sample_no B=sample_name C=sample_id D=section_id E=condition
$data = @(
    @(1,  "Sample_158_A1", "F2", "F2_1", "Flight"),
    @(2,  "Sample_158_B1", "F2", "F2_2", "Flight"),
    @(3,  "Sample_158_C1", "F3", "F3_1", "Flight"),
    @(4,  "Sample_158_D1", "F3", "F3_2", "Flight"),
    @(5,  "Sample_159_A1", "G1", "G1_1", "Ground"),
    @(6,  "Sample_159_B1", "G1", "G1_2", "Ground"),
    @(7,  "Sample_159_C1", "F1", "F1_1", "Flight"),
    @(8,  "Sample_159_D1", "F1", "F1_2", "Flight"),
    @(9,  "Sample_304_A1", "G3", "G3_1", "Ground"),
    @(10, "Sample_304_B1", "G3", "G3_2", "Ground"),
    @(11, "Sample_304_C1", "G2", "G2_1", "Ground"),
    @(12, "Sample_304_D1", "G2", "G2_2", "Ground")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---- drop the old per-row font override on column D --------------------
# (previously some D cells carried a Times New Roman / fontId=1 style)
$ws.Range("A1:E13").ClearFormats()

# ---- column sizing ------------------------------------------------------
# column B (sample_name) keeps its original best-fit width untouched;
# C (sample_id) / D (section_id) are new/resized to fit their longest value
$ws.Columns.Item(3).ColumnWidth = 20.8333333333333
$ws.Columns.Item(4).ColumnWidth = 20.6666666666667

# ---- selection / used range bookkeeping ---------------------------------
$ws.Range("D28").Select()
